# Update ticket/sales-count figures (column F) and mark one listing as
# sold out (column G) across the "展览" (Exhibitions), "演出" (Shows) and
# "全部类型" (All types, the combined roll-up sheet) worksheets.
# Note: "本地生活" (Local life) sheet has no changes in this update.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 33
$ws.Range("F3").Value  = 900
$ws.Range("F4").Value  = 740
$ws.Range("F8").Value  = 738
$ws.Range("F9").Value  = 1106
$ws.Range("F10").Value = 12384
$ws.Range("F20").Value = 504
$ws.Range("F21").Value = 200
$ws.Range("F22").Value = 114
$ws.Range("F23").Value = 320
$ws.Range("F24").Value = 215
$ws.Range("F26").Value = 101
$ws.Range("F27").Value = 109
$ws.Range("F29").Value = 193
$ws.Range("F30").Value = 227
$ws.Range("F31").Value = 1230

# ---- 演出 (Shows) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G3").Value  = "已售罄"
$ws.Range("F6").Value  = 265
$ws.Range("F8").Value  = 114
$ws.Range("F12").Value = 337

# ---- 全部类型 (All types - combined roll-up of the other sheets) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 33
$ws.Range("G4").Value  = "已售罄"
$ws.Range("F5").Value  = 900
$ws.Range("F6").Value  = 740
$ws.Range("F11").Value = 1106
$ws.Range("F12").Value = 12384
$ws.Range("F13").Value = 265
$ws.Range("F22").Value = 504
$ws.Range("F24").Value = 200
$ws.Range("F25").Value = 114
$ws.Range("F26").Value = 114
$ws.Range("F28").Value = 114
$ws.Range("F31").Value = 337
$ws.Range("F32").Value = 320
$ws.Range("F34").Value = 215
$ws.Range("F36").Value = 101
$ws.Range("F37").Value = 109
$ws.Range("F40").Value = 193
$ws.Range("F43").Value = 227
$ws.Range("F44").Value = 1230
